$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.149.69'
$ws.Range('E2').Value = '  -7.81%  '
$ws.Range('D3').Value = '3.643.12'
$ws.Range('E3').Value = '  -7.87%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'575.31"
$ws.Range('E5').Value = '  -5.35%  '
$ws.Range('D6').Value = "'168.44"
$ws.Range('E6').Value = '  -1.58%  '
$ws.Range('D7').Value = '3.629.34'
$ws.Range('E7').Value = '  -8.06%  '
$ws.Range('D8').Value = "'0.613"
$ws.Range('E9').Value = '  +0.14%  '
$ws.Range('D10').Value = "'0.690"
$ws.Range('E10').Value = '  -12.37%  '
$ws.Range('D11').Value = "'0.157"
$ws.Range('E11').Value = '  -13.03%  '
$ws.Range('D12').Value = "'49.92"
$ws.Range('E12').Value = '  -11.74%  '
$ws.Range('D13').Value = "'0.0000281"
$ws.Range('E13').Value = '  -14.18%  '
$ws.Range('D14').Value = "'10.25"
$ws.Range('E14').Value = '  -11.47%  '
$ws.Range('D15').Value = '4.215.80'
$ws.Range('E15').Value = '  -7.87%  '
$ws.Range('D16').Value = '3.672.65'
$ws.Range('E16').Value = '  -7.11%  '
$ws.Range('E17').Value = '  -3.64%  '
$ws.Range('D18').Value = "'19.02"
$ws.Range('E18').Value = '  -10.14%  '
$ws.Range('D19').Value = "'12.63"
$ws.Range('E19').Value = '  -10.17%  '
$ws.Range('D20').Value = "'1.09"
$ws.Range('E20').Value = '  -11.34%  '
$ws.Range('D21').Value = '66.952.81'
$ws.Range('E21').Value = '  -7.85%  '
$ws.Range('D22').Value = "'398.97"
$ws.Range('E22').Value = '  -10.71%  '
$ws.Range('D23').Value = "'4.41"
$ws.Range('E23').Value = '  -9.21%  '
$ws.Range('D24').Value = "'86.42"
$ws.Range('E24').Value = '  -10.12%  '
$ws.Range('E25').Value = '  -11.29%  '
$ws.Range('D26').Value = "'12.50"
$ws.Range('E26').Value = '  -12.32%  '
$ws.Range('D27').Value = "'10.48"
$ws.Range('E27').Value = '  -8.19%  '
$ws.Range('E28').Value = '  +1.40%  '
$ws.Range('D29').Value = "'3.68"
$ws.Range('E29').Value = '  -13.71%  '
$ws.Range('D30').Value = "'9.25"
$ws.Range('E30').Value = '  -11.09%  '
$ws.Range('D31').Value = "'31.95"
$ws.Range('E31').Value = '  -10.85%  '
$ws.Range('D32').Value = "'7.32"
$ws.Range('E32').Value = '  -9.47%  '
$ws.Range('D33').Value = "'12.14"
$ws.Range('E33').Value = '  -12.51%  '
$ws.Range('D34').Value = "'64.10"
$ws.Range('E34').Value = '  -7.04%  '
$ws.Range('E35').Value = '  -11.58%  '
$ws.Range('D36').Value = "'42.29"
$ws.Range('E36').Value = '  -16.27%  '
$ws.Range('D37').Value = "'581.49"
$ws.Range('E37').Value = '  -8.60%  '
$ws.Range('D38').Value = '0.0₃0874'
$ws.Range('E38').Value = '  -13.10%  '
$ws.Range('D39').Value = "'1.00"
$ws.Range('E39').Value = '  +0.10%  '
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('D41').Value = "'0.386"
$ws.Range('E41').Value = '  -10.13%  '
$ws.Range('E42').Value = '  -10.38%  '
$ws.Range('D43').Value = "'2.92"
$ws.Range('E43').Value = '  -15.39%  '
$ws.Range('D44').Value = "'2.58"
$ws.Range('E44').Value = '  -1.67%  '
$ws.Range('D45').Value = "'0.0426"
$ws.Range('E45').Value = '  -11.08%  '
$ws.Range('D46').Value = "'2.79"
$ws.Range('E46').Value = '  -11.82%  '
$ws.Range('D47').Value = "'8.98"
$ws.Range('E47').Value = '  -14.93%  '
$ws.Range('D48').Value = '2.750.55'
$ws.Range('E48').Value = '  -3.03%  '
$ws.Range('E49').Value = '  -11.48%  '
$ws.Range('D50').Value = "'3.12"
$ws.Range('E50').Value = '  -7.82%  '
$ws.Range('D51').Value = "'2.62"
$ws.Range('E51').Value = '  -9.59%  '
